$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the solvent name for row 2: "Nefta" -> "Nafta" Pesada Hidrotratada
$ws.Range("G2").Value = "Nafta Pesada Hidrotratada"

# Widen columns F, H and I to fit their (header) content
$ws.Columns("F").ColumnWidth = 20
$ws.Columns("H").ColumnWidth = 28.3
$ws.Columns("I").ColumnWidth = 23.6

# Move the active selection to G2
$ws.Range("G2").Select() | Out-Null
